$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = "2026-02-21 04:18:15"
    "H2" = "'57%"
    "O2" = "-1.2 °C"
    "E3" = "2026-02-21 04:18:17"
    "H3" = "'63%"
    "M3" = "-0.4 °C 3:52 TU"
    "O3" = "-1.9 °C"
    "E4" = "2026-02-21 04:18:20"
    "K4" = "-0.1 MJ/m2"
    "N4" = "0.9 °C 3:46 TU"
    "O4" = "2.2 °C"
    "E5" = "2026-02-21 04:18:22"
    "H5" = "'65%"
    "M5" = "1.1 °C 3:59 TU"
    "O5" = "-1.2 °C"
    "E6" = "2026-02-21 04:18:24"
    "J6" = "1028.8 hPa"
    "L6" = "13.0 km/h - 339º 3:33 TU"
    "E7" = "2026-02-21 04:18:27"
    "J7" = "1028.3 hPa"
    "O7" = "11.5 °C"
    "E8" = "2026-02-21 04:18:29"
    "J8" = "1028.5 hPa"
    "N8" = "7.0 °C 3:30 TU"
    "O8" = "7.7 °C"
    "E9" = "2026-02-21 04:18:32"
    "N9" = "11.0 °C 3:59 TU"
    "E10" = "2026-02-21 04:18:34"
    "E11" = "2026-02-21 04:18:36"
    "H11" = "'58%"
    "O11" = "4.9 °C"
    "E12" = "2026-02-21 04:18:39"
    "H12" = "'56%"
    "N12" = "9.7 °C 3:44 TU"
    "O12" = "11.9 °C"
    "E13" = "2026-02-21 04:18:41"
    "H13" = "'89%"
    "J13" = "1035.4 hPa"
    "N13" = "-4.6 °C 3:46 TU"
    "O13" = "-3.0 °C"
    "E14" = "2026-02-21 04:18:43"
    "E15" = "2026-02-21 04:18:46"
    "E16" = "2026-02-21 04:18:48"
    "H16" = "'34%"
    "M16" = "0.9 °C 3:43 TU"
    "O16" = "0.0 °C"
    "E17" = "2026-02-21 04:18:50"
    "G17" = "1 cm"
    "N17" = "5.9 °C 3:32 TU"
    "E18" = "2026-02-21 04:18:53"
    "N18" = "0.3 °C 3:31 TU"
    "O18" = "1.1 °C"
    "E19" = "2026-02-21 04:18:55"
    "E20" = "2026-02-21 04:18:57"
    "H20" = "'55%"
    "O20" = "-0.9 °C"
    "E21" = "2026-02-21 04:19:00"
    "H21" = "'74%"
    "J21" = "1032.7 hPa"
    "N21" = "-0.1 °C 3:59 TU"
    "O21" = "1.0 °C"
    "E22" = "2026-02-21 04:19:02"
    "H22" = "'41%"
    "M22" = "-0.5 °C 3:57 TU"
    "O22" = "-1.8 °C"
    "E23" = "2026-02-21 04:19:04"
    "H23" = "'36%"
    "M23" = "1.0 °C 3:59 TU"
    "O23" = "-0.1 °C"
    "E24" = "2026-02-21 04:19:07"
    "N24" = "-0.2 °C 3:39 TU"
    "O24" = "2.2 °C"
    "E25" = "2026-02-21 04:19:09"
    "H25" = "'46%"
    "E26" = "2026-02-21 04:19:11"
    "J26" = "1027.1 hPa"
    "N26" = "5.3 °C 3:35 TU"
    "O26" = "6.5 °C"
    "E27" = "2026-02-21 04:19:14"
    "N27" = "0.4 °C 3:51 TU"
    "E28" = "2026-02-21 04:19:16"
    "J28" = "1030.1 hPa"
    "N28" = "-0.3 °C 3:51 TU"
    "O28" = "0.8 °C"
    "E29" = "2026-02-21 04:19:19"
    "E30" = "2026-02-21 04:19:21"
    "H30" = "'67%"
    "N30" = "9.4 °C 3:45 TU"
    "E31" = "2026-02-21 04:19:23"
    "J31" = "1026.8 hPa"
    "E32" = "2026-02-21 04:19:26"
    "N32" = "-0.2 °C 3:55 TU"
    "O32" = "1.3 °C"
    "E33" = "2026-02-21 04:19:28"
    "H33" = "'77%"
    "J33" = "1033.2 hPa"
    "N33" = "-1.4 °C 3:59 TU"
    "O33" = "-0.3 °C"
    "E34" = "2026-02-21 04:19:31"
    "H34" = "'37%"
    "M34" = "4.0 °C 3:58 TU"
    "E35" = "2026-02-21 04:19:33"
    "J35" = "1031.7 hPa"
    "N35" = "2.6 °C 3:55 TU"
    "E36" = "2026-02-21 04:19:36"
    "N36" = "11.9 °C 3:52 TU"
    "E37" = "2026-02-21 04:19:38"
    "N37" = "-1.7 °C 3:55 TU"
    "E38" = "2026-02-21 04:19:40"
    "L38" = "11.2 km/h - 278º 3:52 TU"
    "O38" = "4.3 °C"
    "E39" = "2026-02-21 04:19:43"
    "H39" = "'39%"
    "M39" = "1.2 °C 3:55 TU"
    "O39" = "0.1 °C"
    "E40" = "2026-02-21 04:19:45"
    "H40" = "'73%"
    "M40" = "4.8 °C 3:31 TU"
    "O40" = "2.5 °C"
    "E41" = "2026-02-21 04:19:48"
    "H41" = "'66%"
    "N41" = "5.1 °C 3:56 TU"
    "O41" = "9.1 °C"
    "E42" = "2026-02-21 04:19:50"
    "E43" = "2026-02-21 04:19:52"
    "N43" = "-0.7 °C 3:46 TU"
    "O43" = "0.6 °C"
    "E44" = "2026-02-21 04:19:54"
    "H44" = "'44%"
    "L44" = "42.5 km/h - 69º 3:30 TU"
    "M44" = "1.7 °C 3:59 TU"
    "O44" = "0.5 °C"
    "E45" = "2026-02-21 04:19:57"
    "H45" = "'88%"
    "J45" = "1035.3 hPa"
    "O45" = "0.5 °C"
    "E46" = "2026-02-21 04:19:59"
    "H46" = "'83%"
    "J46" = "1031.1 hPa"
    "N46" = "4.3 °C 3:58 TU"
    "O46" = "5.7 °C"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
